$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Add a "browser" column (D) to the OpenAccountTest block and populate it,
# and fix the Runmode for the jyoti k / Dollar row from "N" to "Y"
# (removing the now-unused @BeforeSuite toggle scenario).
$ws.Range("D7").Value = "browser"
$ws.Range("D8").Value = "chrome"
$ws.Range("D9").Value = "firefox"

$ws.Range("A9").Value = "Y"

$ws.Range("D7:D9").Select()
